# SPM prep added phase 4
# Appends the "BC_005" subject block (rows 38-49) to Sheet1, mirroring the
# existing per-subject blocks (subj_id in col A, Trial label in col B,
# Start_phase1/Stop_phase1/Start_phase4/Stop_phase4 in C:F, plus a stray
# extra value in G for the ULIFT_L_002 trial row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# subj_id / Trial labels for the new block (rows 38-49)
$labels = @(
    "ULIFT_L_001",
    "ULIFT_L_002",
    "ULIFT_L_003",
    "ULIFT_R_001",
    "ULIFT_R_002",
    "ULIFT_R_003",
    "ULIFT_R_L_001",
    "ULIFT_R_L_002",
    "ULIFT_R_L_003",
    "ULIFT_R_R_001",
    "ULIFT_R_R_002",
    "ULIFT_R_R_003"
)

$startRow = 38
for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "BC_005"
    $ws.Cells.Item($r, 2).Value = $labels[$i]
}

# Numeric measurements (Start_phase1, Stop_phase1, Start_phase4, Stop_phase4)
$ws.Cells.Item(39, 3).Value = 46
$ws.Cells.Item(39, 4).Value = 355
$ws.Cells.Item(39, 5).Value = 1104
$ws.Cells.Item(39, 6).Value = 1402
$ws.Cells.Item(39, 7).Value = 70

$ws.Cells.Item(40, 3).Value = 130
$ws.Cells.Item(40, 4).Value = 406
$ws.Cells.Item(40, 5).Value = 1124
$ws.Cells.Item(40, 6).Value = 1440

$ws.Cells.Item(41, 3).Value = 49
$ws.Cells.Item(41, 4).Value = 399
$ws.Cells.Item(41, 5).Value = 1161
$ws.Cells.Item(41, 6).Value = 1496

$ws.Cells.Item(42, 3).Value = 53
$ws.Cells.Item(42, 4).Value = 393
$ws.Cells.Item(42, 5).Value = 1124
$ws.Cells.Item(42, 6).Value = 1454

$ws.Cells.Item(43, 3).Value = 114
$ws.Cells.Item(43, 4).Value = 425
$ws.Cells.Item(43, 5).Value = 1137
$ws.Cells.Item(43, 6).Value = 1439

# Match the saved cursor/selection state from the edit (I24)
$ws.Range("I24").Select()
